$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.088873982429504
$ws.Range("B1").Value = 2.15859842300415
$ws.Range("C1").Value = 9.295186996459961
$ws.Range("D1").Value = 1.050938725471497
$ws.Range("E1").Value = 1.169913291931152
